$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 2 (weekly refresh: newest reading goes on top, the rest
# of the history shifts down by one row; old row 70 becomes row 71).
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits formatting from the row above (header,
# bold). Clear that and reapply only the date format that column D carries
# throughout the rest of the table.
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new record.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44599
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112038
$ws.Range("G2").Value = "Cebollín baby"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 2400
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = 2450
$ws.Range("N2").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 1225
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = "Hortaliza"
